$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 11, column C: was stored as text "13052054965" -> becomes a real number ---
$ws.Range("C11").Value = 13052054965

# --- 2) Append new conversation/media rows 12-35 ---
# Columns: A=Timestamp, B=Sender, C=Phone, D=Message, E=Media, F=Channel
$phone = 13052054965
$sender = "nnn222111"

$rows = @(
    @{ Row = 12; Ts = "2024-12-03 21:18:38"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-18-38.jpg" }
    @{ Row = 13; Ts = "2024-12-03 21:19:47"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-19-47.jpg" }
    @{ Row = 14; Ts = "2024-12-03 21:20:30"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-20-30.jpg" }
    @{ Row = 15; Ts = "2024-12-03 21:21:03"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-21-03.jpg" }
    @{ Row = 16; Ts = "2024-12-03 21:27:43"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-27-44.jpg" }
    @{ Row = 17; Ts = "2024-12-03 21:30:36"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-30-36.jpg" }
    @{ Row = 18; Ts = "2024-12-03 21:35:14"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-35-14.jpg" }
    @{ Row = 19; Ts = "2024-12-03 21:37:00"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-37-00.jpg" }
    @{ Row = 20; Ts = "2024-12-03 21:47:10"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-47-11.jpg" }
    @{ Row = 21; Ts = "2024-12-03 21:47:55"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-47-55.jpg" }
    @{ Row = 22; Ts = "2024-12-03 21:48:58"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-48-58.jpg" }
    @{ Row = 23; Ts = "2024-12-03 21:50:23"; Msg = "Hello"; Media = $null }
    @{ Row = 24; Ts = "2024-12-03 21:56:55"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-56-55.jpg" }
    @{ Row = 25; Ts = "2024-12-03 21:58:16"; Msg = " ";     Media = "media_files\photo_2024-12-04_02-58-16.jpg" }
    @{ Row = 26; Ts = "2024-12-03 22:02:00"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-02-00.jpg" }
    @{ Row = 27; Ts = "2024-12-03 22:03:10"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-03-11.jpg" }
    @{ Row = 28; Ts = "2024-12-03 22:04:18"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-04-19.jpg" }
    @{ Row = 29; Ts = "2024-12-03 22:04:45"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-04-45.jpg" }
    @{ Row = 30; Ts = "2024-12-03 22:06:39"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-06-39.jpg" }
    @{ Row = 31; Ts = "2024-12-03 22:07:15"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-07-15.jpg" }
    @{ Row = 32; Ts = "2024-12-03 22:08:49"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-08-50.jpg" }
    @{ Row = 33; Ts = "2024-12-03 22:11:11"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-11-12.jpg" }
    @{ Row = 34; Ts = "2024-12-03 22:12:17"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-12-18.jpg" }
    @{ Row = 35; Ts = "2024-12-03 22:12:36"; Msg = " ";     Media = "media_files\photo_2024-12-04_03-12-37.jpg" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Ts
    $ws.Cells.Item($row, 2).Value = $sender

    if ($row -eq 35) {
        # Last added row keeps Phone as literal text, like row 11 originally did.
        $ws.Cells.Item($row, 3).Value = "'13052054965"
    } else {
        $ws.Cells.Item($row, 3).Value = $phone
    }

    $ws.Cells.Item($row, 4).Value = $r.Msg

    if ($r.Media) {
        $ws.Cells.Item($row, 5).Value = $r.Media
    } else {
        # Empty-string media cell (matches the empty inlineStr cells elsewhere).
        $ws.Cells.Item($row, 5).Value = "'"
    }

    # Channel column is always blank for these rows.
    $ws.Cells.Item($row, 6).Value = "'"
}
